$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(8)
$sa = $sh.SmartArt

$node1 = $sa.AllNodes.Item(1)
$node1.TextFrame2.TextRange.Text = "GitHub: https://github.com/LucianoBampa/gerenciador-tarefas"

$node2 = $sa.AllNodes.Item(2)
$node2.TextFrame2.TextRange.Text = "Obrigado!"
